# Automatic update of files.
# - Bump the "Förändrad" (column C) date by 1 day for every data row (2..99)
# - Rewrite HYPERLINK formulas in columns S,T,V,W,X,Y for rows 2..4
#   from folder "Logging_OSTERSUND" to "Logging_2380"
# - Rewrite HYPERLINK formulas in columns S,T,V,W,X,Y for row 5
#   from folder "Logging_BRACKE" to "Logging_2305"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

# Map of row number -> (old folder name, new folder name) for the rows whose
# hyperlink formulas need their project folder renamed.
$folderRenames = @{
    2 = @{ Old = "Logging_OSTERSUND"; New = "Logging_2380" }
    3 = @{ Old = "Logging_OSTERSUND"; New = "Logging_2380" }
    4 = @{ Old = "Logging_OSTERSUND"; New = "Logging_2380" }
    5 = @{ Old = "Logging_BRACKE";    New = "Logging_2305" }
}

for ($row = 2; $row -le $lastRow; $row++) {
    # Bump column C (Förändrad) by one day (45207 -> 45208, etc.)
    $cCell = $ws.Range("C$row")
    $currentValue = $cCell.Value2
    if ($currentValue -ne $null) {
        $cCell.Value2 = $currentValue + 1
    }

    # Update hyperlink formulas that reference the old project folder name.
    if ($folderRenames.ContainsKey($row)) {
        $oldFolder = $folderRenames[$row].Old
        $newFolder = $folderRenames[$row].New

        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Range("$col$row")
            $formula = $cell.Formula
            if ($formula -ne $null -and $formula -ne "") {
                $cell.Formula = $formula.Replace($oldFolder, $newFolder)
            }
        }
    }
}
